# "adding SEs to cookstoves"
# Populate the (previously blank) standard-error column C on the
# wrapper_ready sheet for the carbon_per_cookstove, takeup_control and
# takeup_treatment rows with 0, mirroring the already-populated C5 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("wrapper_ready")

$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 0

# Leave the cursor where the author's saved file shows it.
$ws.Range("H14").Select() | Out-Null
